$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update odds/value cells per the commit diff (only the changed cells are touched).

# Row 5
$ws.Cells.Item(5, 17).Value = 1.98  # Q5: 2 -> 1.98
$ws.Cells.Item(5, 18).Value = 1.88  # R5: 1.85 -> 1.88
$ws.Cells.Item(5, 55).Value = 151  # BC5: 126 -> 151

# Row 6
$ws.Cells.Item(6, 13).Value = 1.06  # M6: 1.07 -> 1.06
$ws.Cells.Item(6, 14).Value = 10  # N6: 9 -> 10
$ws.Cells.Item(6, 17).Value = 2.05  # Q6: 2.1 -> 2.05
$ws.Cells.Item(6, 18).Value = 1.8  # R6: 1.73 -> 1.8
$ws.Cells.Item(6, 19).Value = 1.4  # S6: 1.44 -> 1.4
$ws.Cells.Item(6, 20).Value = 2.75  # T6: 2.63 -> 2.75
$ws.Cells.Item(6, 21).Value = 1.73  # U6: 1.8 -> 1.73
$ws.Cells.Item(6, 22).Value = 2  # V6: 1.91 -> 2
$ws.Cells.Item(6, 29).Value = 9.5  # AC6: 9 -> 9.5
$ws.Cells.Item(6, 34).Value = 8.5  # AH6: 8 -> 8.5
$ws.Cells.Item(6, 44).Value = 67  # AR6: 81 -> 67
$ws.Cells.Item(6, 45).Value = 151  # AS6: 201 -> 151
$ws.Cells.Item(6, 46).Value = 2.75  # AT6: 2.63 -> 2.75

# Row 7
$ws.Cells.Item(7, 8).Value = 3.75  # H7: 3.8 -> 3.75
$ws.Cells.Item(7, 11).Value = 2.25  # K7: 2.3 -> 2.25
$ws.Cells.Item(7, 12).Value = 2.5  # L7: 2.4 -> 2.5
$ws.Cells.Item(7, 15).Value = 1.25  # O7: 1.22 -> 1.25
$ws.Cells.Item(7, 16).Value = 4  # P7: 4.33 -> 4
$ws.Cells.Item(7, 17).Value = 1.8  # Q7: 1.73 -> 1.8
$ws.Cells.Item(7, 18).Value = 2  # R7: 2.1 -> 2
$ws.Cells.Item(7, 19).Value = 1.36  # S7: 1.33 -> 1.36
$ws.Cells.Item(7, 20).Value = 3  # T7: 3.25 -> 3
$ws.Cells.Item(7, 21).Value = 1.73  # U7: 1.67 -> 1.73
$ws.Cells.Item(7, 22).Value = 2  # V7: 2.1 -> 2
$ws.Cells.Item(7, 23).Value = 12  # W7: 13 -> 12
$ws.Cells.Item(7, 29).Value = 12  # AC7: 13 -> 12
$ws.Cells.Item(7, 31).Value = 15  # AE7: 13 -> 15
$ws.Cells.Item(7, 33).Value = 201  # AG7: 151 -> 201
$ws.Cells.Item(7, 34).Value = 8  # AH7: 8.5 -> 8
$ws.Cells.Item(7, 38).Value = 15  # AL7: 13 -> 15
$ws.Cells.Item(7, 46).Value = 3  # AT7: 3.25 -> 3
$ws.Cells.Item(7, 52).Value = 34  # AZ7: 29 -> 34

# Row 8
$ws.Cells.Item(8, 7).Value = 2.05  # G8: 2 -> 2.05
$ws.Cells.Item(8, 8).Value = 3.4  # H8: 3.3 -> 3.4
$ws.Cells.Item(8, 9).Value = 3.5  # I8: 3.6 -> 3.5
$ws.Cells.Item(8, 13).Value = 1.05  # M8: 1.06 -> 1.05
$ws.Cells.Item(8, 14).Value = 11  # N8: 10 -> 11
$ws.Cells.Item(8, 15).Value = 1.29  # O8: 1.3 -> 1.29
$ws.Cells.Item(8, 16).Value = 3.75  # P8: 3.5 -> 3.75
$ws.Cells.Item(8, 17).Value = 1.93  # Q8: 2.03 -> 1.93
$ws.Cells.Item(8, 18).Value = 1.93  # R8: 1.83 -> 1.93
$ws.Cells.Item(8, 19).Value = 1.4  # S8: 1.37 -> 1.4
$ws.Cells.Item(8, 21).Value = 1.73  # U8: 1.8 -> 1.73
$ws.Cells.Item(8, 22).Value = 2  # V8: 1.91 -> 2
$ws.Cells.Item(8, 23).Value = 8  # W8: 7.5 -> 8
$ws.Cells.Item(8, 24).Value = 10  # X8: 9.5 -> 10
$ws.Cells.Item(8, 29).Value = 11  # AC8: 10 -> 11
$ws.Cells.Item(8, 31).Value = 13  # AE8: 15 -> 13
$ws.Cells.Item(8, 32).Value = 41  # AF8: 51 -> 41
$ws.Cells.Item(8, 33).Value = 201  # AG8: 251 -> 201
$ws.Cells.Item(8, 36).Value = 12  # AJ8: 13 -> 12
$ws.Cells.Item(8, 47).Value = 7.5  # AU8: 8 -> 7.5
$ws.Cells.Item(8, 50).Value = 19  # AX8: 21 -> 19
$ws.Cells.Item(8, 51).Value = 26  # AY8: 29 -> 26
$ws.Cells.Item(8, 54).Value = 151  # BB8: 201 -> 151

# Row 9
$ws.Cells.Item(9, 19).Value = 1.33  # S9: 1.3 -> 1.33

# Row 10
$ws.Cells.Item(10, 7).Value = 1.91  # G10: 1.85 -> 1.91
$ws.Cells.Item(10, 8).Value = 3.5  # H10: 3.6 -> 3.5
$ws.Cells.Item(10, 9).Value = 3.9  # I10: 4 -> 3.9
$ws.Cells.Item(10, 19).Value = 1.33  # S10: 1.3 -> 1.33
$ws.Cells.Item(10, 23).Value = 8.5  # W10: 9 -> 8.5
$ws.Cells.Item(10, 28).Value = 23  # AB10: 21 -> 23
$ws.Cells.Item(10, 45).Value = 126  # AS10: 101 -> 126
$ws.Cells.Item(10, 48).Value = 51  # AV10: 41 -> 51
$ws.Cells.Item(10, 49).Value = 5.5  # AW10: 6 -> 5.5
$ws.Cells.Item(10, 50).Value = 19  # AX10: 21 -> 19

# Row 11
$ws.Cells.Item(11, 7).Value = 1.91  # G11: 1.95 -> 1.91
$ws.Cells.Item(11, 8).Value = 3.6  # H11: 3.5 -> 3.6
$ws.Cells.Item(11, 9).Value = 3.7  # I11: 3.5 -> 3.7
$ws.Cells.Item(11, 10).Value = 2.5  # J11: 2.6 -> 2.5
$ws.Cells.Item(11, 19).Value = 1.36  # S11: 1.33 -> 1.36
$ws.Cells.Item(11, 21).Value = 1.7  # U11: 1.67 -> 1.7
$ws.Cells.Item(11, 22).Value = 2.05  # V11: 2.1 -> 2.05
$ws.Cells.Item(11, 35).Value = 21  # AI11: 19 -> 21
$ws.Cells.Item(11, 36).Value = 13  # AJ11: 12 -> 13
$ws.Cells.Item(11, 38).Value = 29  # AL11: 26 -> 29
$ws.Cells.Item(11, 41).Value = 10  # AO11: 11 -> 10
$ws.Cells.Item(11, 50).Value = 21  # AX11: 19 -> 21
$ws.Cells.Item(11, 52).Value = 67  # AZ11: 51 -> 67

# Row 12
$ws.Cells.Item(12, 7).Value = 2.55  # G12: 2.45 -> 2.55
$ws.Cells.Item(12, 9).Value = 2.7  # I12: 2.8 -> 2.7
$ws.Cells.Item(12, 10).Value = 3.2  # J12: 3.1 -> 3.2
$ws.Cells.Item(12, 17).Value = 1.98  # Q12: 1.95 -> 1.98
$ws.Cells.Item(12, 18).Value = 1.88  # R12: 1.9 -> 1.88
$ws.Cells.Item(12, 26).Value = 26  # Z12: 23 -> 26
$ws.Cells.Item(12, 27).Value = 21  # AA12: 19 -> 21
$ws.Cells.Item(12, 35).Value = 13  # AI12: 15 -> 13
$ws.Cells.Item(12, 36).Value = 10  # AJ12: 11 -> 10
$ws.Cells.Item(12, 41).Value = 15  # AO12: 13 -> 15
